$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 833
$ws.Range("J17").Value = 837.4888999999999
$ws.Range("L17").Value = 2512.4667
$ws.Range("N17").Value = -2848.4667
$ws.Range("H58").Value = 17858918
$ws.Range("I58").Value = 35714430
$ws.Range("J58").Value = 3404.8572
$ws.Range("K58").Value = 107143290
$ws.Range("L58").Value = 10214.5716
$ws.Range("M58").Value = -107143140
$ws.Range("N58").Value = -10514.5716
$ws.Range("H112").Value = 3539.4
$ws.Range("J112").Value = 3654.889
$ws.Range("L112").Value = 10964.667
$ws.Range("N112").Value = -13180.667
$ws.Range("H132").Value = 1993.2632
$ws.Range("I132").Value = 1591.625
$ws.Range("K132").Value = 4774.875
$ws.Range("M132").Value = -2244.875
$ws.Range("H137").Value = 5232.9688
$ws.Range("I137").Value = 6685.3477
$ws.Range("J137").Value = 1521.3334
$ws.Range("K137").Value = 20056.0431
$ws.Range("L137").Value = 4564.0002
$ws.Range("M137").Value = -17506.0431
$ws.Range("N137").Value = -9664.0002
$ws.Range("H138").Value = 3999.4194
$ws.Range("J138").Value = 4758.905
$ws.Range("L138").Value = 14276.715
$ws.Range("N138").Value = -24556.715

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2360.2036
$ws.Range("I32").Value = 2316.4614
$ws.Range("J32").Value = 3497.5
$ws.Range("K32").Value = 2316.4614
$ws.Range("L32").Value = 3497.5
$ws.Range("M32").Value = -2029.4614
$ws.Range("N32").Value = -4071.5
$ws.Range("H61").Value = 3320.611
$ws.Range("I61").Value = 3320.611
$ws.Range("K61").Value = 3320.611
$ws.Range("M61").Value = -3108.611
$ws.Range("H74").Value = 2164.077
$ws.Range("I74").Value = 2041.1765
$ws.Range("J74").Value = 2999.8
$ws.Range("K74").Value = 2041.1765
$ws.Range("L74").Value = 2999.8
$ws.Range("M74").Value = -1167.1765
$ws.Range("N74").Value = -4747.8
$ws.Range("H77").Value = 2164.077
$ws.Range("I77").Value = 2041.1765
$ws.Range("J77").Value = 2999.8
$ws.Range("K77").Value = 10205.8825
$ws.Range("L77").Value = 14999
$ws.Range("M77").Value = -5837.8825
$ws.Range("N77").Value = -23735
$ws.Range("H97").Value = 1038.6666
$ws.Range("I97").Value = 1052.5
$ws.Range("J97").Value = 1011
$ws.Range("K97").Value = 1052.5
$ws.Range("L97").Value = 1011
$ws.Range("M97").Value = -556.5
$ws.Range("N97").Value = -2003
$ws.Range("H102").Value = 4383.6113
$ws.Range("I102").Value = 2223.4614
$ws.Range("K102").Value = 2223.4614
$ws.Range("M102").Value = -601.4614000000001
$ws.Range("H109").Value = 130188.5
$ws.Range("J109").Value = 130188.5
$ws.Range("L109").Value = 130188.5
$ws.Range("N109").Value = -132962.5
$ws.Range("H122").Value = 6947741.5
$ws.Range("I122").Value = 10103869
$ws.Range("J122").Value = 4259.8
$ws.Range("K122").Value = 30311607
$ws.Range("L122").Value = 12779.4
$ws.Range("M122").Value = -30309157
$ws.Range("N122").Value = -17679.4
$ws.Range("H132").Value = 4100
$ws.Range("I132").Value = 4100
$ws.Range("K132").Value = 12300
$ws.Range("M132").Value = -9770
$ws.Range("H133").Value = 70499.75
$ws.Range("J133").Value = 70499.75
$ws.Range("L133").Value = 70499.75
$ws.Range("N133").Value = -75559.75
$ws.Range("H134").Value = 49999
$ws.Range("J134").Value = 49999
$ws.Range("L134").Value = 49999
$ws.Range("N134").Value = -60139
$ws.Range("H135").Value = 65000
$ws.Range("J135").Value = 65000
$ws.Range("L135").Value = 65000
$ws.Range("N135").Value = -75140
$ws.Range("H136").Value = 3320.611
$ws.Range("I136").Value = 3320.611
$ws.Range("K136").Value = 9961.832999999999
$ws.Range("M136").Value = -7411.832999999999

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1864.2424
$ws.Range("I94").Value = 1668.5714
$ws.Range("J94").Value = 2960
$ws.Range("K94").Value = 1668.5714
$ws.Range("L94").Value = 2960
$ws.Range("M94").Value = -1217.5714
$ws.Range("N94").Value = -3862
$ws.Range("H99").Value = 3532.3076
$ws.Range("I99").Value = 1835
$ws.Range("K99").Value = 1835
$ws.Range("M99").Value = -337
$ws.Range("H107").Value = 5361
$ws.Range("I107").Value = 854.4545000000001
$ws.Range("K107").Value = 854.4545000000001
$ws.Range("M107").Value = 1065.5455
$ws.Range("H134").Value = 4229.7856
$ws.Range("I134").Value = 4309
$ws.Range("K134").Value = 12927
$ws.Range("M134").Value = -10392

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3104.9412
$ws.Range("I31").Value = 1730.8096
$ws.Range("J31").Value = 5324.6924
$ws.Range("K31").Value = 1730.8096
$ws.Range("L31").Value = 5324.6924
$ws.Range("M31").Value = -1435.8096
$ws.Range("N31").Value = -5914.6924
$ws.Range("H34").Value = 3104.9412
$ws.Range("I34").Value = 1730.8096
$ws.Range("J34").Value = 5324.6924
$ws.Range("K34").Value = 1730.8096
$ws.Range("L34").Value = 5324.6924
$ws.Range("M34").Value = -1528.8096
$ws.Range("N34").Value = -5728.6924
$ws.Range("H99").Value = 1834.0526
$ws.Range("I99").Value = 1801.1875
$ws.Range("K99").Value = 1801.1875
$ws.Range("M99").Value = -303.1875
$ws.Range("H107").Value = 1426.6666
$ws.Range("I107").Value = 1442.625
$ws.Range("J107").Value = 1299
$ws.Range("K107").Value = 1442.625
$ws.Range("L107").Value = 1299
$ws.Range("M107").Value = 477.375
$ws.Range("N107").Value = -5139
$ws.Range("H122").Value = 4503.6665
$ws.Range("I122").Value = 4505.5
$ws.Range("K122").Value = 13516.5
$ws.Range("M122").Value = -11066.5
$ws.Range("H126").Value = 1834.0526
$ws.Range("I126").Value = 1801.1875
$ws.Range("K126").Value = 5403.5625
$ws.Range("M126").Value = -2933.5625
$ws.Range("H134").Value = 2203.9092
$ws.Range("I134").Value = 2299.3
$ws.Range("K134").Value = 6897.900000000001
$ws.Range("M134").Value = -4362.900000000001

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2261.158
$ws.Range("I68").Value = 2017.8889
$ws.Range("J68").Value = 2480.1
$ws.Range("K68").Value = 6053.6667
$ws.Range("L68").Value = 7440.299999999999
$ws.Range("M68").Value = -5242.6667
$ws.Range("N68").Value = -9062.299999999999
$ws.Range("H71").Value = 2261.158
$ws.Range("I71").Value = 2017.8889
$ws.Range("J71").Value = 2480.1
$ws.Range("K71").Value = 18161.0001
$ws.Range("L71").Value = 22320.9
$ws.Range("M71").Value = -14105.0001
$ws.Range("N71").Value = -30432.9
$ws.Range("H86").Value = 1958.25
$ws.Range("I86").Value = 1946
$ws.Range("J86").Value = 1970.5
$ws.Range("K86").Value = 5838
$ws.Range("L86").Value = 5911.5
$ws.Range("M86").Value = -4652
$ws.Range("N86").Value = -8283.5
$ws.Range("H89").Value = 1958.25
$ws.Range("I89").Value = 1946
$ws.Range("J89").Value = 1970.5
$ws.Range("K89").Value = 17514
$ws.Range("L89").Value = 17734.5
$ws.Range("M89").Value = -11586
$ws.Range("N89").Value = -29590.5
$ws.Range("H113").Value = 1428.3572
$ws.Range("I113").Value = 1157.9
$ws.Range("K113").Value = 3473.7
$ws.Range("M113").Value = -1303.7

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 45462468
$ws.Range("I113").Value = 200001600
$ws.Range("J113").Value = 9783.117
$ws.Range("K113").Value = 200001600
$ws.Range("L113").Value = 9783.117
$ws.Range("M113").Value = -199999430
$ws.Range("N113").Value = -14123.117
$ws.Range("H122").Value = 13853.7
$ws.Range("I122").Value = 16092.167
$ws.Range("K122").Value = 48276.501
$ws.Range("M122").Value = -45826.501

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3612.111
$ws.Range("I7").Value = 3654
$ws.Range("K7").Value = 3654
$ws.Range("M7").Value = -3542
$ws.Range("H40").Value = 5949.9688
$ws.Range("I40").Value = 3276.4614
$ws.Range("K40").Value = 3276.4614
$ws.Range("M40").Value = -3140.4614
$ws.Range("H93").Value = 5623.1333
$ws.Range("I93").Value = 4534.7
$ws.Range("K93").Value = 4534.7
$ws.Range("M93").Value = -3286.7
$ws.Range("H126").Value = 3612.111
$ws.Range("I126").Value = 3654
$ws.Range("K126").Value = 10962
$ws.Range("M126").Value = -8492

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 10877.777
$ws.Range("I62").Value = 9975
$ws.Range("J62").Value = 11600
$ws.Range("K62").Value = 9975
$ws.Range("L62").Value = 11600
$ws.Range("M62").Value = -9351
$ws.Range("N62").Value = -12848
$ws.Range("H65").Value = 10877.777
$ws.Range("I65").Value = 9975
$ws.Range("J65").Value = 11600
$ws.Range("K65").Value = 49875
$ws.Range("L65").Value = 58000
$ws.Range("M65").Value = -46755
$ws.Range("N65").Value = -64240
$ws.Range("H86").Value = 49990
$ws.Range("J86").Value = 49990
$ws.Range("L86").Value = 49990
$ws.Range("N86").Value = -52236
$ws.Range("H89").Value = 49990
$ws.Range("J89").Value = 49990
$ws.Range("L89").Value = 249950
$ws.Range("N89").Value = -261182
$ws.Range("H122").Value = 9960.4
$ws.Range("I122").Value = 10934.5
$ws.Range("J122").Value = 9311
$ws.Range("K122").Value = 32803.5
$ws.Range("L122").Value = 27933
$ws.Range("M122").Value = -30353.5
$ws.Range("N122").Value = -32833
$ws.Range("H126").Value = 2635.2104
$ws.Range("I126").Value = 2756.8235
$ws.Range("K126").Value = 8270.470499999999
$ws.Range("M126").Value = -5800.470499999999
$ws.Range("H132").Value = 8555.478999999999
$ws.Range("I132").Value = 6614.385
$ws.Range("K132").Value = 19843.155
$ws.Range("M132").Value = -17313.155
$ws.Range("H136").Value = 3406.7856
$ws.Range("I136").Value = 3217.7273
$ws.Range("J136").Value = 4100
$ws.Range("K136").Value = 9653.1819
$ws.Range("L136").Value = 12300
$ws.Range("M136").Value = -7103.1819
$ws.Range("N136").Value = -17400
